$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sat after ", 2016)" (near the
#    Bootstrap citation). Removing it causes Word to renumber the remaining
#    bookmark ids (the two "_Toc..." bookmarks shift down by one), matching
#    the diff automatically.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Rewrite the "Apache" paragraph:
#    "Como servidor web local de pruebas se ha decidido por Apache web
#    server, con el fin de coincidir con el servidor de producción. "
#    becomes
#    "Como servidor web local se ha decidido usar Apache web server."
#    with a new "_GoBack" bookmark inserted between " usar" and
#    " Apache web server".
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*servidor web local de pruebas*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $rng.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
    $startPos = $rng.Start

    $newText = "Como servidor web local se ha decidido usar Apache web server."
    $rng.Text = $newText
    $endPos = $rng.End

    $pUsar   = $startPos + $newText.IndexOf(" usar")
    $pApache = $startPos + $newText.IndexOf(" Apache")
    $pDot    = $startPos + $newText.LastIndexOf(".")

    # Insert bookmarks from right (end of string) to left so the earlier
    # offsets stay valid, and the run before the final "." is not left with
    # a dangling xml:space="preserve" artifact.
    $d.Bookmarks.Add("ZZTMP_SPLIT2", $d.Range($pDot, $pDot)) | Out-Null
    $d.Bookmarks.Add("_GoBack", $d.Range($pApache, $pApache)) | Out-Null
    $d.Bookmarks.Add("ZZTMP_SPLIT1", $d.Range($pUsar, $pUsar)) | Out-Null

    # Drop the helper bookmarks - the run split they forced survives.
    $d.Bookmarks("ZZTMP_SPLIT1").Delete()
    $d.Bookmarks("ZZTMP_SPLIT2").Delete()
}

# ---------------------------------------------------------------------------
# 3. Delete the now-orphaned "3.1.7.1 Características" and
#    "3.1.8 Plataforma de producción" heading paragraphs entirely.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$p1 = $null
$p2 = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "*3.1.7.1*Caracter*") {
        $p1 = $p
    } elseif ($t -like "*3.1.8*Plataforma de producci*") {
        $p2 = $p
    }
}

if (($p1 -ne $null) -and ($p2 -ne $null)) {
    $delRange = $d.Range($p1.Range.Start, $p2.Range.End)
    $delRange.Delete()
}
